# Updates loading_percent values for rows 2-25, columns C,D,E,F,G,I,J,M,N,O
# (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colC = @(3.322466478192388, 3.288864571732341, 3.267733310853226, 3.259000105531092, 3.25754269774724, 3.267616021194663, 3.31098634460686, 3.391919211534769, 3.448673818267319, 3.473865022106858, 3.483310855958623, 3.481280739393291, 3.474644030109244, 3.470566574683885, 3.447014608836152, 3.432403261528401, 3.42394040535601, 3.421065034799943, 3.433964774493416, 3.476595957432745, 3.503911241338503, 3.489383676493509, 3.433259009199467, 3.37048608360743)
$colD = @(8.638992300073987, 8.667074326772095, 8.685219072085934, 8.692840653631105, 8.694119965609975, 8.685320937707097, 8.648488060525034, 8.583391775605556, 8.539876830622928, 8.521009016541985, 8.513997034363864, 8.515501289618186, 8.520429476906404, 8.523465420018709, 8.541128504267565, 8.552201392324312, 8.558657541243958, 8.560858494432603, 8.551013631225363, 8.518978347321356, 8.498815556507505, 8.509506150227686, 8.551550337026843, 8.600242274839282)
$colE = @(12.37201559825996, 12.44140398867836, 12.48605657273915, 12.50476935314296, 12.50790783619372, 12.48630684656521, 12.39551662312596, 12.23365421745352, 12.12449385427441, 12.07693143210543, 12.05922043696315, 12.06302150275741, 12.07546833630959, 12.08313139008072, 12.12764420895655, 12.15548699162147, 12.17169871603598, 12.17722165593352, 12.15250267127151, 12.07180427269934, 12.02081037539214, 12.04786737927213, 12.15385124527181, 12.27572100434)
$colF = @(32.54574313732715, 32.39441187383002, 32.31167324943948, 32.28053947857153, 32.27552633451081, 32.31124288312613, 32.49146658585119, 32.92443636391938, 33.28913270729051, 33.46467596889456, 33.5324901458415, 33.5178263374936, 33.47022850903759, 33.44124650538136, 33.27785062955766, 33.18004966951332, 33.12470760199335, 33.10612747098825, 33.19036684806937, 33.48417317084535, 33.68397803564594, 33.57664220700282, 33.18569969140896, 32.7989712263597)
$colG = @(3.630240520644479, 3.63340281407815, 3.635447463541185, 3.636306663532887, 3.636450905380848, 3.635458945664814, 3.631309559281708, 3.623985627140736, 3.619094539838489, 3.616974573440489, 3.616186803235762, 3.616355797205544, 3.616909462676703, 3.617250551630891, 3.619235190428475, 3.620479537385735, 3.62120514233507, 3.621452520853632, 3.620346051628926, 3.616746430883649, 3.614481350449505, 3.615682290314427, 3.620406368719337, 3.625880511114031)
$colI = @(23.00355010586011, 22.86973139493713, 22.79426592581971, 22.76522257176308, 22.76050385322129, 22.79386728523808, 22.95603463646812, 23.32597497721318, 23.62759102433816, 23.77083575868803, 23.82590657713705, 23.8140099769848, 23.77535013266915, 23.7517763296801, 23.61834767214825, 23.53801134185441, 23.49237492946359, 23.47702250311425, 23.5465044856102, 23.78668335406598, 23.94845510654388, 23.86168921769576, 23.54266301895881, 23.22051614098257)
$colJ = @(9.340452987757853, 9.384523245209513, 9.412859956720116, 9.42472960781188, 9.426720044437333, 9.413018729058557, 9.355383894142086, 9.252451319735201, 9.182911322196514, 9.152582823624545, 9.141284896313294, 9.143709813063939, 9.151649597544392, 9.156537246464982, 9.184919553534346, 9.202664913278291, 9.212994518734602, 9.216513089665543, 9.200763171987674, 9.149312427577906, 9.116774877377667, 9.134041488506856, 9.201622551804917, 9.279224040940614)
$colM = @(23.73462125797119, 23.05034516553914, 22.62023740240928, 22.44269394160919, 22.41308297316443, 22.61785186381432, 23.50088461193904, 25.14411438863251, 26.28619384068628, 26.78955017544916, 26.97768265906496, 26.93727739186719, 26.80507841199399, 26.7237758493987, 26.25295946316044, 25.95987214038857, 25.78978098658817, 25.73193568060229, 25.99122973299017, 26.84397682379355, 27.38679562866824, 27.09845566865926, 25.97705791042704, 24.71027626824297)
$colN = @(17.41703623623913, 17.14195899497876, 16.97280910835753, 16.90389865283502, 16.89245974081672, 16.97187956881292, 17.32229099186326, 18.00389334178022, 18.49665465511209, 18.71815517954719, 18.80157582485754, 18.78363116039325, 18.72502781504861, 18.68906993595304, 18.48211912701639, 18.35442641974386, 18.28073490852963, 18.25574429473078, 18.36804554107394, 18.74225401133652, 18.98412629244543, 18.85530473900339, 18.36188920049834, 17.82058340548246)
$colO = @(24.87948638320081, 24.74858482806647, 24.6755062490762, 24.64758296052244, 24.64305905736912, 24.67512212060014, 24.83285168165171, 25.19891841578698, 25.50076058273719, 25.64479880603449, 25.70027102262187, 25.68828346262399, 25.6493441723044, 25.62561238717235, 25.49147941436932, 25.41088911543603, 25.36517155015168, 25.3498027445098, 25.41940256289688, 25.66075672315419, 25.82388541667477, 25.7363409565462, 25.4155517150365, 25.0939766227154)

$columnData = @{
    3 = $colC
    4 = $colD
    5 = $colE
    6 = $colF
    7 = $colG
    9 = $colI
    10 = $colJ
    13 = $colM
    14 = $colN
    15 = $colO
}

foreach ($colIndex in $columnData.Keys) {
    $values = $columnData[$colIndex]
    $row = 2
    foreach ($v in $values) {
        $ws.Cells.Item($row, $colIndex).Value = $v
        $row = $row + 1
    }
}
